# Generate Report for Handback
# Adds a new handback entry (29f631d0-13d9-4325-bd0a-8b9626ddc855.md) as
# row 4 to the "Overview", "zh-cn" and "de-de" worksheets/tables.

$wb = $excel.ActiveWorkbook

$fileGuid   = "29f631d0-13d9-4325-bd0a-8b9626ddc855"
$fileName   = "$fileGuid.md"
$pathName   = "e2e\$fileGuid.md"
$ext        = ".md"
$statusSync = "Handed back: in sync with en-US"

$srcCommit  = "45d8bf2db12b74ff31aa99e92a25ca1b1a628d53"
$zhCommit   = "e2880f13dd5dfe07cd125ca6084a52b8a7b9f319"
$deCommit   = "375e7d6a290b0f98ffe497a6e71088b2754be420"

$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$fileGuid.md"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhCommit/e2e/$fileGuid.md"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$deCommit/e2e/$fileGuid.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = "2016-11-09 06:49:49"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcUrl, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = $fileName
$wsZhCn.Range("B4").Value = $ext
$wsZhCn.Range("C4").Value = $statusSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = "$fileGuid.d96602369b5deaf6d97f6c87ebe3abe6a746005c.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-11-09 06:49:35"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = $fileName
$wsZhCn.Range("J4").Value = "$fileGuid.d96602369b5deaf6d97f6c87ebe3abe6a746005c.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-11-09 06:50:27"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), $zhUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = $fileName
$wsDeDe.Range("B4").Value = $ext
$wsDeDe.Range("C4").Value = $statusSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = "$fileGuid.d96602369b5deaf6d97f6c87ebe3abe6a746005c.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-11-09 06:49:49"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = $fileName
$wsDeDe.Range("J4").Value = "$fileGuid.d96602369b5deaf6d97f6c87ebe3abe6a746005c.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-11-09 06:50:45"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), $deUrl, "", "", $fileName) | Out-Null
